$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the remaining log entry's Student ID and Log Date
$ws.Range("A2").Value = "201255"
$ws.Range("C2").Value = "20/09/2025"

# The other logged excuses (rows 3-9) are gone from the export - remove them
$ws.Range("A3:F9").EntireRow.Delete()
